$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "shopping with Ansel Adams"

# B2 holds a plain text date string ("2013-06-15"); force Text format so
# Excel doesn't reinterpret it as a date serial number.
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "2013-06-15"

$ws.Range("E2").Value = "Heavy Duty Concrete Shirt"
$ws.Range("F2").Value = "Turnips"

# H2 holds the text "2.0" (not the number 2); force Text format so Excel
# doesn't reinterpret it as a number.
$ws.Range("H2").NumberFormat = "@"
$ws.Range("H2").Value = "2.0"
